$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as text so values like "1.027" or
# "5.480" are not reinterpreted as numbers and lose formatting/precision.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.392.84'
$ws.Range("E2").Value = '  +3.40%  '
$ws.Range("D3").Value = '1.841.45'
$ws.Range("E3").Value = '  +3.66%  '
$ws.Range("D4").Value = '1.027'
$ws.Range("E4").Value = '  +2.92%  '
$ws.Range("D5").Value = '319.29'
$ws.Range("E5").Value = '  +2.13%  '
$ws.Range("D6").Value = '1.023'
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("D7").Value = '0.4353'
$ws.Range("E7").Value = '  +1.59%  '
$ws.Range("D8").Value = '0.3721'
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("D9").Value = '0.07345'
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").Value = '0.8764'
$ws.Range("E10").Value = '  +3.03%  '
$ws.Range("D11").Value = '21.42'
$ws.Range("E11").Value = '  +4.07%  '
$ws.Range("D12").Value = '2.000.96'
$ws.Range("E12").Value = '  +13.15%  '
$ws.Range("D13").Value = '5.480'
$ws.Range("E13").Value = '  +4.03%  '
$ws.Range("D14").Value = '6.680'
$ws.Range("E14").Value = '  +3.39%  '
$ws.Range("D15").Value = '0.07155'
$ws.Range("E15").Value = '  +3.69%  '
$ws.Range("D16").Value = '82.13'
$ws.Range("E16").Value = '  +3.87%  '
$ws.Range("D17").Value = '1.028'
$ws.Range("E17").Value = '  +3.11%  '
$ws.Range("D18").Value = '0.000008994'
$ws.Range("E18").Value = '  +3.29%  '
$ws.Range("D19").Value = '1.023'
$ws.Range("E19").Value = '  +2.57%  '
$ws.Range("D20").Value = '15.42'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").Value = '27.430.10'
$ws.Range("E21").Value = '  +3.59%  '
$ws.Range("D22").Value = '5.251'
$ws.Range("E22").Value = '  +2.45%  '
$ws.Range("D23").Value = '11.14'
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = '2.201.53'
$ws.Range("E24").Value = '  +10.92%  '
$ws.Range("D25").Value = '156.82'
$ws.Range("E25").Value = '  +2.86%  '
$ws.Range("D26").Value = '1.909'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = '18.53'
$ws.Range("E27").Value = '  +2.44%  '
$ws.Range("D28").Value = '5.285'
$ws.Range("E28").Value = '  +3.44%  '
$ws.Range("D29").Value = '1.925'
$ws.Range("E29").Value = '  +5.92%  '
$ws.Range("D30").Value = '115.49'
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = '0.09021'
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").Value = '  +5.80%  '
$ws.Range("D33").Value = '0.7602'
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("D34").Value = '4.466'
$ws.Range("E34").Value = '  +2.90%  '
$ws.Range("D35").Value = '2.856'
$ws.Range("E35").Value = '  +4.35%  '
$ws.Range("D36").Value = '1.025'
$ws.Range("E36").Value = '  +2.80%  '
$ws.Range("D37").Value = '1.150'
$ws.Range("E37").Value = '  +3.32%  '
$ws.Range("D38").Value = '0.01956'
$ws.Range("E38").Value = '  +3.43%  '
$ws.Range("D39").Value = '0.05250'
$ws.Range("E39").Value = '  +1.54%  '
$ws.Range("D40").Value = '0.5159'
$ws.Range("E40").Value = '  +4.29%  '
$ws.Range("D41").Value = '2.800'
$ws.Range("E41").Value = '  +7.27%  '
$ws.Range("D42").Value = '0.1662'
$ws.Range("E42").Value = '  +2.90%  '
$ws.Range("D43").Value = '6.535'
$ws.Range("E43").Value = '  +3.32%  '
$ws.Range("D44").Value = '8.482'
$ws.Range("E44").Value = '  +5.29%  '
$ws.Range("D45").Value = '108.13'
$ws.Range("E45").Value = '  +2.67%  '
$ws.Range("D46").Value = '10.56'
$ws.Range("E46").Value = '  +3.63%  '
$ws.Range("D47").Value = '1.026'
$ws.Range("E47").Value = '  +2.88%  '
$ws.Range("D48").Value = '0.4633'
$ws.Range("E48").Value = '  +2.77%  '
$ws.Range("D49").Value = '1.669'
$ws.Range("E49").Value = '  +2.47%  '
$ws.Range("D50").Value = '1.889'
$ws.Range("E50").Value = '  +8.06%  '
$ws.Range("D51").Value = '0.06289'
$ws.Range("E51").Value = '  +1.37%  '
